$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1070
$ws.Range("J97").Value = 1070
$ws.Range("L97").Value = 3210
$ws.Range("N97").Value = -4202
$ws.Range("H98").Value = 2409.15
$ws.Range("I98").Value = 1565.5333
$ws.Range("J98").Value = 4940
$ws.Range("K98").Value = 1565.5333
$ws.Range("L98").Value = 4940
$ws.Range("M98").Value = -67.53330000000005
$ws.Range("N98").Value = -7936
$ws.Range("H122").Value = 2409.15
$ws.Range("I122").Value = 1565.5333
$ws.Range("J122").Value = 4940
$ws.Range("K122").Value = 4696.5999
$ws.Range("L122").Value = 14820
$ws.Range("M122").Value = -2246.5999
$ws.Range("N122").Value = -19720
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 68171.836
$ws.Range("I23").Value = 80006
$ws.Range("K23").Value = 80006
$ws.Range("M23").Value = -79747
$ws.Range("H32").Value = 405721.28
$ws.Range("I32").Value = 461942.38
$ws.Range("J32").Value = 12173.429
$ws.Range("K32").Value = 461942.38
$ws.Range("L32").Value = 12173.429
$ws.Range("M32").Value = -461655.38
$ws.Range("N32").Value = -12747.429
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 2359.238
$ws.Range("I45").Value = 1853.6
$ws.Range("K45").Value = 1853.6
$ws.Range("M45").Value = -1476.6
$ws.Range("H63").Value = 1507.25
$ws.Range("I63").Value = 1505.5
$ws.Range("J63").Value = 1512.5
$ws.Range("K63").Value = 1505.5
$ws.Range("L63").Value = 1512.5
$ws.Range("M63").Value = -819.5
$ws.Range("N63").Value = -2884.5
$ws.Range("H66").Value = 1507.25
$ws.Range("I66").Value = 1505.5
$ws.Range("J66").Value = 1512.5
$ws.Range("K66").Value = 7527.5
$ws.Range("L66").Value = 7562.5
$ws.Range("M66").Value = -4095.5
$ws.Range("N66").Value = -14426.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1542.909
$ws.Range("I58").Value = 1117.25
$ws.Range("J58").Value = 1786.1428
$ws.Range("K58").Value = 1117.25
$ws.Range("L58").Value = 1786.1428
$ws.Range("M58").Value = -914.25
$ws.Range("N58").Value = -2192.1428
$ws.Range("H60").Value = 29350
$ws.Range("J60").Value = 29350
$ws.Range("L60").Value = 29350
$ws.Range("N60").Value = -30372
$ws.Range("H132").Value = 29417338
$ws.Range("I132").Value = 55563628
$ws.Range("J132").Value = 2760.75
$ws.Range("K132").Value = 166690884
$ws.Range("L132").Value = 8282.25
$ws.Range("M132").Value = -166688354
$ws.Range("N132").Value = -13342.25
$ws.Range("H134").Value = 3478.0557
$ws.Range("I134").Value = 4240.6665
$ws.Range("J134").Value = 1952.8334
$ws.Range("K134").Value = 12721.9995
$ws.Range("L134").Value = 5858.5002
$ws.Range("M134").Value = -10186.9995
$ws.Range("N134").Value = -10928.5002
$ws.Range("H136").Value = 1542.909
$ws.Range("I136").Value = 1117.25
$ws.Range("J136").Value = 1786.1428
$ws.Range("K136").Value = 3351.75
$ws.Range("L136").Value = 5358.428400000001
$ws.Range("M136").Value = -801.75
$ws.Range("N136").Value = -10458.4284
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 761.619
$ws.Range("J34").Value = 920.13336
$ws.Range("L34").Value = 2760.40008
$ws.Range("N34").Value = -2928.40008
$ws.Range("H39").Value = 3436.3635
$ws.Range("J39").Value = 3700
$ws.Range("L39").Value = 11100
$ws.Range("N39").Value = -11688
$ws.Range("H48").Value = 2413.3333
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 2413.3333
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 7239.999899999999
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -7739.999899999999
$ws.Range("H55").Value = 2654.9
$ws.Range("J55").Value = 3571.4285
$ws.Range("L55").Value = 10714.2855
$ws.Range("N55").Value = -11068.2855
$ws.Range("H68").Value = 1616.8182
$ws.Range("I68").Value = 753.7778
$ws.Range("J68").Value = 5500.5
$ws.Range("K68").Value = 2261.3334
$ws.Range("L68").Value = 16501.5
$ws.Range("M68").Value = -1450.3334
$ws.Range("N68").Value = -18123.5
$ws.Range("H71").Value = 1616.8182
$ws.Range("I71").Value = 753.7778
$ws.Range("J71").Value = 5500.5
$ws.Range("K71").Value = 6784.000199999999
$ws.Range("L71").Value = 49504.5
$ws.Range("M71").Value = -2728.000199999999
$ws.Range("N71").Value = -57616.5
$ws.Range("H92").Value = 903
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 903
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2709
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5205
$ws.Range("H131").Value = 752.7455
$ws.Range("I131").Value = 290
$ws.Range("J131").Value = 942.5897
$ws.Range("K131").Value = 870
$ws.Range("L131").Value = 2827.7691
$ws.Range("M131").Value = 4170
$ws.Range("N131").Value = -12907.7691
$ws.Range("H137").Value = 6007418
$ws.Range("I137").Value = 101394.55
$ws.Range("J137").Value = 16835128
$ws.Range("K137").Value = 304183.65
$ws.Range("L137").Value = 50505384
$ws.Range("M137").Value = -299083.65
$ws.Range("N137").Value = -50515584
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2589.3333
$ws.Range("I80").Value = 2141.111
$ws.Range("J80").Value = 3037.5557
$ws.Range("K80").Value = 2141.111
$ws.Range("L80").Value = 3037.5557
$ws.Range("M80").Value = -1143.111
$ws.Range("N80").Value = -5033.5557
$ws.Range("H83").Value = 2589.3333
$ws.Range("I83").Value = 2141.111
$ws.Range("J83").Value = 3037.5557
$ws.Range("K83").Value = 10705.555
$ws.Range("L83").Value = 15187.7785
$ws.Range("M83").Value = -5713.555
$ws.Range("N83").Value = -25171.7785
$ws.Range("H122").Value = 1874
$ws.Range("I122").Value = 1369.7
$ws.Range("J122").Value = 2434.3333
$ws.Range("K122").Value = 4109.1
$ws.Range("L122").Value = 7302.999899999999
$ws.Range("M122").Value = -1659.1
$ws.Range("N122").Value = -12202.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 244.45454
$ws.Range("I55").Value = 205.88889
$ws.Range("J55").Value = 271.15384
$ws.Range("K55").Value = 205.88889
$ws.Range("L55").Value = 271.15384
$ws.Range("M55").Value = -32.88889
$ws.Range("N55").Value = -617.1538399999999
$ws.Range("H63").Value = 11000
$ws.Range("H66").Value = 11000
$ws.Range("H68").Value = 1648.1482
$ws.Range("I68").Value = 1705.8823
$ws.Range("J68").Value = 1550
$ws.Range("K68").Value = 1705.8823
$ws.Range("L68").Value = 1550
$ws.Range("M68").Value = -956.8823
$ws.Range("N68").Value = -3048
$ws.Range("H71").Value = 1648.1482
$ws.Range("I71").Value = 1705.8823
$ws.Range("J71").Value = 1550
$ws.Range("K71").Value = 8529.4115
$ws.Range("L71").Value = 7750
$ws.Range("M71").Value = -4785.4115
$ws.Range("N71").Value = -15238
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3182.5
$ws.Range("I62").Value = 2992
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 2992
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2368
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3182.5
$ws.Range("I65").Value = 2992
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 14960
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -11840
$ws.Range("N65").Value = -23740
$ws.Range("H122").Value = 40001320
$ws.Range("I122").Value = 100000600
$ws.Range("J122").Value = 1798.3334
$ws.Range("K122").Value = 300001800
$ws.Range("L122").Value = 5395.0002
$ws.Range("M122").Value = -299999350
$ws.Range("N122").Value = -10295.0002

Write-Host "Applied 202 cell updates across 7 sheets"
